# Applies the Zeromus_Profits market-price/profit recalculation update
# produced by the scheduled Sheets runner: refreshed currentAveragePrice
# (NQ/HQ) values pulled in, with LevePrice/LeveProfit columns recomputed
# to match for each affected leve row.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")

# Row 98
$ws.Range("H98").Value = 422.14285
$ws.Range("I98").Value = 369.16666
$ws.Range("J98").Value = 740
$ws.Range("K98").Value = 369.16666
$ws.Range("L98").Value = 740
$ws.Range("M98").Value = 1128.83334
$ws.Range("N98").Value = -3736

# Row 122
$ws.Range("H122").Value = 422.14285
$ws.Range("I122").Value = 369.16666
$ws.Range("J122").Value = 740
$ws.Range("K122").Value = 1107.49998
$ws.Range("L122").Value = 2220
$ws.Range("M122").Value = 1342.50002
$ws.Range("N122").Value = -7120

# Row 133
$ws.Range("H133").Value = 46166.668
$ws.Range("J133").Value = 46166.668
$ws.Range("L133").Value = 46166.668
$ws.Range("N133").Value = -56286.668

# Row 137
$ws.Range("H137").Value = 2010.5385
$ws.Range("I137").Value = 1648.8182
$ws.Range("K137").Value = 4946.4546
$ws.Range("M137").Value = -2396.4546

$ws = $wb.Worksheets.Item("ARM")

# Row 61
$ws.Range("H61").Value = 1851.125
$ws.Range("I61").Value = 1739.381
$ws.Range("K61").Value = 1739.381
$ws.Range("M61").Value = -1527.381

# Row 62
$ws.Range("H62").Value = 15000
$ws.Range("J62").Value = 15000
$ws.Range("L62").Value = 15000
$ws.Range("N62").Value = -16248

# Row 65
$ws.Range("H65").Value = 15000
$ws.Range("J65").Value = 15000
$ws.Range("L65").Value = 45000
$ws.Range("N65").Value = -51240

# Row 68
$ws.Range("H68").Value = 18330
$ws.Range("J68").Value = 18330
$ws.Range("L68").Value = 18330
$ws.Range("N68").Value = -19952

# Row 71
$ws.Range("H71").Value = 18330
$ws.Range("J71").Value = 18330
$ws.Range("L71").Value = 54990
$ws.Range("N71").Value = -63102

# Row 74
$ws.Range("H74").Value = 7410411
$ws.Range("I74").Value = 8698199
$ws.Range("J74").Value = 5632
$ws.Range("K74").Value = 8698199
$ws.Range("L74").Value = 5632
$ws.Range("M74").Value = -8697325
$ws.Range("N74").Value = -7380

# Row 75
$ws.Range("H75").Value = 40000
$ws.Range("J75").Value = 40000
$ws.Range("L75").Value = 40000
$ws.Range("N75").Value = -41748

# Row 77
$ws.Range("H77").Value = 7410411
$ws.Range("I77").Value = 8698199
$ws.Range("J77").Value = 5632
$ws.Range("K77").Value = 43490995
$ws.Range("L77").Value = 28160
$ws.Range("M77").Value = -43486627
$ws.Range("N77").Value = -36896

# Row 78
$ws.Range("H78").Value = 40000
$ws.Range("J78").Value = 40000
$ws.Range("L78").Value = 120000
$ws.Range("N78").Value = -128736

# Row 80
$ws.Range("H80").Value = 36666.668
$ws.Range("I80").Value = 0
$ws.Range("J80").Value = 36666.668
$ws.Range("K80").Value = 0
$ws.Range("L80").Value = 36666.668
$ws.Range("M80").ClearContents()
$ws.Range("N80").Value = -38662.668

# Row 81
$ws.Range("H81").Value = 37700
$ws.Range("J81").Value = 37700
$ws.Range("L81").Value = 37700
$ws.Range("N81").Value = -39696

# Row 82
$ws.Range("H82").Value = 23120.666
$ws.Range("I82").Value = 0
$ws.Range("J82").Value = 23120.666
$ws.Range("K82").Value = 0
$ws.Range("L82").Value = 23120.666
$ws.Range("M82").ClearContents()
$ws.Range("N82").Value = -23842.666

# Row 83
$ws.Range("H83").Value = 36666.668
$ws.Range("I83").Value = 0
$ws.Range("J83").Value = 36666.668
$ws.Range("K83").Value = 0
$ws.Range("L83").Value = 110000.004
$ws.Range("M83").ClearContents()
$ws.Range("N83").Value = -119984.004

# Row 84
$ws.Range("H84").Value = 37700
$ws.Range("J84").Value = 37700
$ws.Range("L84").Value = 113100
$ws.Range("N84").Value = -123084

# Row 85
$ws.Range("H85").Value = 23120.666
$ws.Range("I85").Value = 0
$ws.Range("J85").Value = 23120.666
$ws.Range("K85").Value = 0
$ws.Range("L85").Value = 23120.666
$ws.Range("M85").ClearContents()
$ws.Range("N85").Value = -25616.666

# Row 86
$ws.Range("H86").Value = 23095
$ws.Range("I86").Value = 23095
$ws.Range("K86").Value = 23095
$ws.Range("M86").Value = -21909

# Row 89
$ws.Range("H89").Value = 23095
$ws.Range("I89").Value = 23095
$ws.Range("K89").Value = 69285
$ws.Range("M89").Value = -63357

# Row 136
$ws.Range("H136").Value = 1851.125
$ws.Range("I136").Value = 1739.381
$ws.Range("K136").Value = 5218.143
$ws.Range("M136").Value = -2668.143

$ws = $wb.Worksheets.Item("BSM")

# Row 134
$ws.Range("H134").Value = 3627.182
$ws.Range("I134").Value = 3000
$ws.Range("J134").Value = 3862.375
$ws.Range("K134").Value = 9000
$ws.Range("L134").Value = 11587.125
$ws.Range("M134").Value = -6465
$ws.Range("N134").Value = -16657.125

$ws = $wb.Worksheets.Item("CUL")

# Row 11
$ws.Range("H11").Value = 161.8125
$ws.Range("I11").Value = 179.875
$ws.Range("J11").Value = 143.75
$ws.Range("K11").Value = 539.625
$ws.Range("L11").Value = 431.25
$ws.Range("M11").Value = -399.625
$ws.Range("N11").Value = -711.25

$ws = $wb.Worksheets.Item("GSM")

# Row 102
$ws.Range("H102").Value = 1813.0526
$ws.Range("I102").Value = 1913.4286
$ws.Range("K102").Value = 1913.4286
$ws.Range("M102").Value = -291.4286

# Row 122
$ws.Range("H122").Value = 1988.9166
$ws.Range("I122").Value = 2149.8333
$ws.Range("K122").Value = 6449.499899999999
$ws.Range("M122").Value = -3999.499899999999

# Row 132
$ws.Range("H132").Value = 5999
$ws.Range("I132").Value = 0
$ws.Range("J132").Value = 5999
$ws.Range("K132").Value = 0
$ws.Range("L132").Value = 17997
$ws.Range("M132").ClearContents()
$ws.Range("N132").Value = -23057

$ws = $wb.Worksheets.Item("LTW")

# Row 7
$ws.Range("H7").Value = 1759.875
$ws.Range("I7").Value = 1327.4546
$ws.Range("J7").Value = 2711.2
$ws.Range("K7").Value = 1327.4546
$ws.Range("L7").Value = 2711.2
$ws.Range("M7").Value = -1215.4546
$ws.Range("N7").Value = -2935.2

# Row 40
$ws.Range("H40").Value = 1780.1333
$ws.Range("I40").Value = 1493.3334
$ws.Range("J40").Value = 1851.8334
$ws.Range("K40").Value = 1493.3334
$ws.Range("L40").Value = 1851.8334
$ws.Range("M40").Value = -1357.3334
$ws.Range("N40").Value = -2123.8334

# Row 46
$ws.Range("H46").Value = 919.5
$ws.Range("I46").Value = 333.66666
$ws.Range("J46").Value = 1054.6923
$ws.Range("K46").Value = 333.66666
$ws.Range("L46").Value = 1054.6923
$ws.Range("M46").Value = -145.66666
$ws.Range("N46").Value = -1430.6923

# Row 126
$ws.Range("H126").Value = 1759.875
$ws.Range("I126").Value = 1327.4546
$ws.Range("J126").Value = 2711.2
$ws.Range("K126").Value = 3982.3638
$ws.Range("L126").Value = 8133.599999999999
$ws.Range("M126").Value = -1512.3638
$ws.Range("N126").Value = -13073.6

$ws = $wb.Worksheets.Item("WVR")

# Row 122
$ws.Range("H122").Value = 90910760
$ws.Range("I122").Value = 250001120
$ws.Range("J122").Value = 1986.4286
$ws.Range("K122").Value = 750003360
$ws.Range("L122").Value = 5959.2858
$ws.Range("M122").Value = -750000910
$ws.Range("N122").Value = -10859.2858
